$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 249334
$ws.Range("E3").Value = 1036486441

$ws.Range("C91").Value = 151202
$ws.Range("E91").Value = 482926544

$ws.Range("C92").Value = 409279
$ws.Range("E92").Value = 1597133234

$ws.Range("C93").Value = 209653
$ws.Range("E93").Value = 1309991539

$ws.Range("C94").Value = 94229
$ws.Range("E94").Value = 918982287

$ws.Range("C95").Value = 50801
$ws.Range("E95").Value = 934231255

$ws.Range("C96").Value = 17322
$ws.Range("E96").Value = 797242656

$ws.Range("C104").Value = 135300
$ws.Range("E104").Value = 272650477

$ws.Range("C167").Value = 12220
$ws.Range("E167").Value = 105794259
